$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2120343839541547
$ws.Range("C2").Value = 0.5186246418338109
$ws.Range("J2").Value = 0.0171919770773639
$ws.Range("P2").Value = 0.1547277936962751
$ws.Range("S2").Value = 0.09742120343839542
$ws.Range("B3").Value = 0.02162162162162162
$ws.Range("C3").Value = 0.02162162162162162
$ws.Range("J3").Value = 0.02162162162162162
$ws.Range("P3").Value = 0.7297297297297297
$ws.Range("S3").Value = 0.2054054054054054
$ws.Range("B6").Value = 0.08294930875576037
$ws.Range("E6").Value = 0.004608294930875576
$ws.Range("F6").Value = 0.04608294930875576
$ws.Range("J6").Value = 0.2764976958525346
$ws.Range("O6").Value = 0.03225806451612903
$ws.Range("Q6").Value = 0.1658986175115207
$ws.Range("R6").Value = 0.09677419354838709
$ws.Range("S6").Value = 0.2949308755760369
$ws.Range("B7").Value = 0.1180555555555556
$ws.Range("D7").Value = 0.02083333333333333
$ws.Range("F7").Value = 0.0625
$ws.Range("J7").Value = 0.1805555555555556
$ws.Range("O7").Value = 0.0625
$ws.Range("Q7").Value = 0.1458333333333333
$ws.Range("R7").Value = 0.09027777777777778
$ws.Range("S7").Value = 0.3194444444444444
$ws.Range("B8").Value = 0.09740259740259741
$ws.Range("D8").Value = 0.01948051948051948
$ws.Range("F8").Value = 0.04112554112554113
$ws.Range("J8").Value = 0.1038961038961039
$ws.Range("O8").Value = 0.02813852813852814
$ws.Range("Q8").Value = 0.20995670995671
$ws.Range("R8").Value = 0.1103896103896104
$ws.Range("S8").Value = 0.3896103896103896
$ws.Range("B9").Value = 0.07391304347826087
$ws.Range("D9").Value = 0.03478260869565217
$ws.Range("F9").Value = 0.09565217391304348
$ws.Range("J9").Value = 0.1260869565217391
$ws.Range("O9").Value = 0.01304347826086956
$ws.Range("Q9").Value = 0.2391304347826087
$ws.Range("R9").Value = 0.0782608695652174
$ws.Range("S9").Value = 0.3391304347826087
$ws.Range("B10").Value = 0.1294027565084227
$ws.Range("D10").Value = 0.02526799387442573
$ws.Range("E10").Value = 0.003062787136294028
$ws.Range("F10").Value = 0.06891271056661562
$ws.Range("J10").Value = 0.1079632465543645
$ws.Range("O10").Value = 0.01607963246554365
$ws.Range("Q10").Value = 0.2166921898928025
$ws.Range("R10").Value = 0.09341500765696784
$ws.Range("S10").Value = 0.3392036753445635
$ws.Range("G11").Value = 0.1137440758293839
$ws.Range("J11").Value = 0.1327014218009479
$ws.Range("K11").Value = 0.1990521327014218
$ws.Range("L11").Value = 0.5450236966824644
$ws.Range("S11").Value = 0.009478672985781991
$ws.Range("G12").Value = 0.717741935483871
$ws.Range("J12").Value = 0.1935483870967742
$ws.Range("K12").Value = 0.01612903225806452
$ws.Range("L12").Value = 0.03225806451612903
$ws.Range("S12").Value = 0.04032258064516129
$ws.Range("F15").Value = 0.02192982456140351
$ws.Range("H15").Value = 0.1403508771929824
$ws.Range("I15").Value = 0.07894736842105263
$ws.Range("J15").Value = 0.4035087719298245
$ws.Range("K15").Value = 0.008771929824561403
$ws.Range("M15").Value = 0.0131578947368421
$ws.Range("O15").Value = 0.04385964912280702
$ws.Range("S15").Value = 0.2894736842105263
$ws.Range("F16").Value = 0.0136986301369863
$ws.Range("H16").Value = 0.1506849315068493
$ws.Range("I16").Value = 0.0547945205479452
$ws.Range("J16").Value = 0.4292237442922374
$ws.Range("K16").Value = 0.091324200913242
$ws.Range("M16").Value = 0.0365296803652968
$ws.Range("O16").Value = 0.0639269406392694
$ws.Range("S16").Value = 0.1598173515981735
$ws.Range("F17").Value = 0.01642710472279261
$ws.Range("H17").Value = 0.1971252566735113
$ws.Range("I17").Value = 0.1375770020533881
$ws.Range("J17").Value = 0.431211498973306
$ws.Range("K17").Value = 0.06776180698151951
$ws.Range("M17").Value = 0.01848049281314168
$ws.Range("N17").Value = 0.002053388090349076
$ws.Range("O17").Value = 0.05749486652977413
$ws.Range("S17").Value = 0.07186858316221766
$ws.Range("F18").Value = 0.01345291479820628
$ws.Range("H18").Value = 0.1973094170403587
$ws.Range("I18").Value = 0.1076233183856502
$ws.Range("J18").Value = 0.4170403587443946
$ws.Range("K18").Value = 0.06278026905829596
$ws.Range("M18").Value = 0.01345291479820628
$ws.Range("O18").Value = 0.07174887892376682
$ws.Range("S18").Value = 0.1165919282511211
$ws.Range("F19").Value = 0.01808785529715762
$ws.Range("H19").Value = 0.227390180878553
$ws.Range("I19").Value = 0.0921619293712317
$ws.Range("J19").Value = 0.4013781223083548
$ws.Range("K19").Value = 0.08440999138673558
$ws.Range("M19").Value = 0.01981050818260121
$ws.Range("O19").Value = 0.06890611541774333
$ws.Range("S19").Value = 0.08785529715762273
